$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) values for the rows whose price changed
$priceUpdates = @{
    2 = "260.46"
    3 = "22.74"
    4 = "6.191"
    5 = "0.06075"
    6 = "6.735"
    7 = "3.452"
    8 = "1.352"
    9 = "0.7979"
    10 = "0.1587"
    11 = "0.08032"
    12 = "0.03359"
    13 = "0.03058"
    14 = "0.09318"
    15 = "3.851"
    16 = "0.001702"
    17 = "0.04863"
    18 = "0.0006140"
    19 = "0.006222"
    20 = "0.001093"
    21 = "0.003395"
    23 = "3.712"
    24 = "2.255"
    25 = "0.3362"
    26 = "0.1262"
    27 = "0.0003202"
    40 = "0.04575"
    41 = "0.007090"
    42 = "0.003900"
    43 = "0.1116"
    44 = "0.01006"
    45 = "0.002970"
    46 = "0.00005946"
    48 = "0.7000"
    49 = "0.07075"
    50 = "0.00002100"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item([int]$row, 4)   # column D = Price
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

# Column G ("Hora") goes from "0" to "2" for every data row (2-51)
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = Hora
    $cell.NumberFormat = "@"
    $cell.Value = "2"
    $cell.Style = "Normal"
}
